$wb = $excel.ActiveWorkbook

# Add a new worksheet "NoHeaders" at the end of the workbook, containing the
# same data rows as the "Another" sheet but without the header row.
$count = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($count))
$ws.Name = "NoHeaders"

$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = $true
$ws.Range("D1").Value = 44715
$ws.Range("D1").NumberFormat = "d-mmm"

$ws.Range("A2").Value = "b"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = $false
$ws.Range("D2").Value = 44693
$ws.Range("D2").NumberFormat = "d-mmm"

$ws.Range("A3").Value = "c"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = 44607
$ws.Range("D3").NumberFormat = "d-mmm"

$ws.Range("A1:D3").Select() | Out-Null

# Update the selection on the "Another" sheet to cover the data rows,
# with D4 (bottom-right corner) as the active cell within that selection.
$another = $wb.Worksheets.Item("Another")
$another.Activate() | Out-Null
$another.Range("D4").Activate() | Out-Null
$another.Range("A2:D4").Select() | Out-Null

# Make the newly added sheet the active tab, matching the target workbook
# view (activeTab points at the 3rd sheet / NoHeaders).
$ws.Activate() | Out-Null
